$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-12-19"

$ws.Range("I1").Value = "2022 (through 12-19)"
$ws.Range("I13").Value = 81
$ws.Range("I14").Value = 1598
